$d = $word.ActiveDocument

# Plain find & replace, scoped to a given range. Re-fetches the scope range
# text live each time (caller passes the Range object fresh or re-derived),
# safe to call multiple times in sequence against the same logical paragraph.
function Replace-Text($scopeRange, $old, $new) {
    $scopeRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# "Safe" replace: finds $old within $scopeRange, inserts $new immediately
# *after* the found text (preserving any markers such as commentRangeStart
# that sit right before the found text, and the formatting of any run that
# sits right after it), then clears the original found text. This avoids
# two quirks of this engine's plain Find+Replace:
#   1) a comment anchor sitting at the exact start of the found text gets
#      shifted to after the replacement text
#   2) replacement text placed right after a hyperlink run borrows the
#      hyperlink's run formatting instead of keeping the original run's
#      formatting
function Safe-Replace-Text($scopeRange, $old, $new) {
    $find = $scopeRange.Duplicate
    $ok = $find.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return
    }
    $insPoint = $find.Duplicate
    $insPoint.Collapse(0)
    $insPoint.InsertBefore($new)
    $find.Text = ""
}

# --- Paragraph 1: language switcher line (hyperlinked "English" + list) ---
Replace-Text $d.Paragraphs(1).Range "English" "Inglês"
Safe-Replace-Text $d.Paragraphs(1).Range " / Portuguese / French / Thai / Vietnamese / Spanish" " / Português / Francês / Tailandês / Vietnamita / Espanhol"

# --- Paragraph 3: standalone "English" heading ---
Replace-Text $d.Paragraphs(3).Range "English" "Inglês"

# --- Paragraph 5: "Brief:" ---
Replace-Text $d.Paragraphs(5).Range "Brief" "Resumo"

# --- Paragraph 6: brief description ---
Replace-Text $d.Paragraphs(6).Range "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io" "Um e-mail enviado a parceiros no país-alvo que confirmaram presença, mas não submeteram os seus documentos até à data limite. Vamos revogar os seus convites. Será enviado através do customer.io"

# --- Paragraph 8: "Target audience:" ---
Replace-Text $d.Paragraphs(8).Range "Target audience" "Público-alvo"

# --- Paragraph 9: target audience description ---
Replace-Text $d.Paragraphs(9).Range "Invited partners who didn’t submit their documents on time" "Parceiros convidados que não submeteram os seus documentos atempadamente"

# --- Paragraph 12: Subject line (drop the trailing " registration" run) ---
Replace-Text $d.Paragraphs(12).Range "Subject line" "Linha de assunto"
Replace-Text $d.Paragraphs(12).Range ": Your " ": A sua inscrição na "
Replace-Text $d.Paragraphs(12).Range "[EVENT NAME]" "[NOME DO EVENTO]"
Replace-Text $d.Paragraphs(12).Range " registration" ""

# --- Paragraph 14: "We didn't receive your documents on time" heading ---
Replace-Text $d.Paragraphs(14).Range "We didn’t receive your documents on time" "Não recebemos os seus documentos a tempo"

# --- Paragraph 16: "Hi [PARTNER NAME], " ---
Replace-Text $d.Paragraphs(16).Range "Hi " "Olá "
Replace-Text $d.Paragraphs(16).Range "[PARTNER NAME]" "[NOME DO PARCEIRO]"

# --- Paragraph 17: deadline paragraph ---
Replace-Text $d.Paragraphs(17).Range "We didn’t receive your documents by the deadline (" "Não recebemos os seus documentos dentro do prazo ("
Replace-Text $d.Paragraphs(17).Range "). Unfortunately, we’re unable to proceed with your registration for the " "). Infelizmente, não é possível prosseguir com o seu registo para a "
Replace-Text $d.Paragraphs(17).Range "[EVENT NAME]" "[NOME DO EVENTO]"

# --- Paragraph 19: live chat / WhatsApp contact (comment 1 starts here, and
#     a hyperlink immediately follows the first run) ---
Safe-Replace-Text $d.Paragraphs(19).Range "If you have any questions, please contact us via " "Para mais informações, contacte-nos através de "
Safe-Replace-Text $d.Paragraphs(19).Range " or " " ou "

# --- Paragraph 20: country manager contact (comment 1 ends here) ---
Replace-Text $d.Paragraphs(20).Range "If you have any questions, please contact your country manager, " "Para mais questões, pode também contactar o seus gestor de parcerias "
Replace-Text $d.Paragraphs(20).Range ", at " ", em "
Replace-Text $d.Paragraphs(20).Range " or " " ou "

# --- Comments: "choose either one" -> "escolha um de cada vez" (both) ---
foreach ($c in $d.Comments) {
    Replace-Text $c.Range "choose either one" "escolha um de cada vez"
}
